$wb = $excel.ActiveWorkbook

# --- Fix typo in sheet name: "Unit Dimains" -> "Unit Domains" ---
$wb.Worksheets.Item("Unit Dimains").Name = "Unit Domains"

# --- Data Dictionary: clarify a few column descriptions ---
$dd = $wb.Worksheets.Item("Data Dictionary")
$dd.Range("B2").Value = "Name of the wqs parameter."
$dd.Range("B3").Value = "Name of the wqs parameter as reported in ECHO. Restricted to only those listed in the Parameter Domains tab"
$dd.Range("B6").Value = "Unit of parameter as reported in ECHO. Restricted to only those listed in the Unit Domains tab."

# --- Move the active selection: Data Dictionary ends up with B9 selected ---
$dd.Range("B9").Select() | Out-Null

# --- Make the WQS sheet the active/selected tab, with E8 selected ---
$wqs = $wb.Worksheets.Item("WQS")
$wqs.Activate() | Out-Null
$wqs.Range("E8").Select() | Out-Null
